$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = [double]"1.303490370430959E-14"
$ws.Range("E3").Value = [double]"1.303490370430959E-14"

# Row 4
$ws.Range("D4").Value = [double]"0.9734303091231553"
$ws.Range("E4").Value = [double]"0.9734303091231553"

# Row 5
$ws.Range("D5").Value = [double]"3.231907912200274E-17"
$ws.Range("E5").Value = [double]"3.231907912200274E-17"

# Row 6
$ws.Range("D6").Value = [double]"2.575603430795212E-38"
$ws.Range("E6").Value = [double]"2.575603430795212E-38"

# Row 7
$ws.Range("D7").Value = [double]"1.497982658358462E-06"
$ws.Range("E7").Value = [double]"0.9999985020173416"

# Row 9
$ws.Range("D9").Value = [double]"0.9999999989252772"
$ws.Range("E9").Value = [double]"1.074722755234347E-09"

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"4.534341690877145E-06"
$ws.Range("E10").Value = [double]"0.9999954656583091"

# Row 11
$ws.Range("D11").Value = [double]"8.038394568237695E-10"
$ws.Range("E11").Value = [double]"0.9999999991961606"
$ws.Range("F11").Value = [double]"9.421274185180664"
$ws.Range("G11").Value = [double]"0.5"
